# Update of the GL results
# - SVM sheet: add STDEV row (row 24) under the existing AVERAGE row (row 23)
# - GL_MV sheet: drop the now-stray "STD" label row (row 25)
# - GL_adaptive sheet: replace the stray "STD" label row (row 25) with a
#   proper STD row (row 24) that carries STDEV formulas across B:F

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "SVM": add row 24 with STDEV formulas under the AVERAGE row (23)
# ---------------------------------------------------------------------
$wsSVM = $wb.Worksheets.Item("SVM")
$wsSVM.Range("A24").Formula = "=STDEV(A3:A12)"
$wsSVM.Range("B24").Formula = "=STDEV(B3:B12)"

# ---------------------------------------------------------------------
# Sheet "GL_MV": remove the stray "STD" label row (row 25)
# ---------------------------------------------------------------------
$wsGLMV = $wb.Worksheets.Item("GL_MV")
$wsGLMV.Range("A25").ClearContents()

# ---------------------------------------------------------------------
# Sheet "GL_adaptive": remove the stray "STD" label row (row 25) and
# add a real STD row (row 24) with STDEV formulas across B:F
# ---------------------------------------------------------------------
$wsGLA = $wb.Worksheets.Item("GL_adaptive")
$wsGLA.Range("A25").ClearContents()

$wsGLA.Range("A24").Value = "STD"
$wsGLA.Range("B24").Formula = "=STDEV(B3:B12)"
$wsGLA.Range("C24:F24").Formula = "=STDEV(C3:C12)"

# Restore selection state to match the committed workbook
$wsSVM.Range("F21").Select()
$wsGLMV.Range("A24").Select()
$wsGLA.Range("B24:F24").Select()
$wsGLA.Activate()
